$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 224; this shifts the existing rows 224..268 down to 225..269
$ws.Rows.Item(224).Insert()

# Populate the new row 224 with the new weekly price record (same Mercado/Region/Categoria as the rest of the block)
$ws.Cells.Item(224, 1).Value2 = 4
$ws.Cells.Item(224, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(224, 3).Value2 = "Los Lagos"
$ws.Cells.Item(224, 4).Value2 = 44722
$ws.Cells.Item(224, 5).Value2 = 10
$ws.Cells.Item(224, 6).Value2 = 100112021
$ws.Cells.Item(224, 7).Value2 = "Ají"
$ws.Cells.Item(224, 8).Value2 = "Inferno"
$ws.Cells.Item(224, 9).Value2 = "Primera"
$ws.Cells.Item(224, 10).Value2 = 120
$ws.Cells.Item(224, 11).Value2 = 24000
$ws.Cells.Item(224, 12).Value2 = 24000
$ws.Cells.Item(224, 13).Value2 = 24000
$ws.Cells.Item(224, 14).Value2 = "$/caja 12 kilos"
$ws.Cells.Item(224, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(224, 16).Value2 = 2000
$ws.Cells.Item(224, 17).Value2 = 12
$ws.Cells.Item(224, 18).Value2 = "Hortaliza"
